# Add "game 12" results as a new row 13 on the glory worksheet, mirroring
# the formula pattern already present in row 12 (fill-down style).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("glory")

# --- raw game inputs (columns A-G) ---
$ws.Range("A13").Value2 = 12
$ws.Range("B13").Value2 = 22
$ws.Range("C13").Value2 = 15
$ws.Range("D13").Value2 = 25
$ws.Range("E13").Value2 = 17
$ws.Range("F13").Value2 = 16
$ws.Range("G13").Value2 = 15

# --- running totals (columns H-L), same shape as row 12's formulas ---
$ws.Range("H13").Formula = '=SUM(C$2:C13)'
$ws.Range("I13").Formula = '=SUM(D$2:D13)'
$ws.Range("J13").Formula = '=SUM(E$2:E13)'
$ws.Range("K13").Formula = '=SUM(F$2:F13)'
$ws.Range("L13").Formula = '=SUM(G$2:G13)'

# --- diff from the game leader (columns M-Q) ---
$ws.Range("M13").Formula = '=H13-MAX(H13:L13)'
$ws.Range("N13").Formula = '=I13-MAX(H13:L13)'
$ws.Range("O13").Formula = '=J13-MAX(H13:L13)'
$ws.Range("P13").Formula = '=K13-MAX(H13:L13)'
$ws.Range("Q13").Formula = '=L13-MAX(H13:L13)'

# --- per-game rank (columns R-V) ---
$ws.Range("R13").Formula = '=RANK(C13,$C13:$G13)'
$ws.Range("S13").Formula = '=RANK(D13,$C13:$G13)'
$ws.Range("T13").Formula = '=RANK(E13,$C13:$G13)'
$ws.Range("U13").Formula = '=RANK(F13,$C13:$G13)'
$ws.Range("V13").Formula = '=RANK(G13,$C13:$G13)'

# --- normalized score vs target (columns W-AA) ---
$ws.Range("W13").Formula = '=C13/$B13'
$ws.Range("X13").Formula = '=D13/$B13'
$ws.Range("Y13").Formula = '=E13/$B13'
$ws.Range("Z13").Formula = '=F13/$B13'
$ws.Range("AA13").Formula = '=G13/$B13'

# --- total glory for the game (column AB) ---
$ws.Range("AB13").Formula = '=SUM(C13:G13)'

# Match the author's final selection on the sheet.
$ws.Range("G13").Select()
